# "Generate Report for Archive"
#
# Two changes, matching the canonical OOXML diff:
#   1. Every "Ready for handoff" status cell becomes "In Translation"
#      (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4 — all 12 cells reference
#      the same shared string, so every usage must be updated together).
#   2. The "Status"/locale columns narrow from ~17.22 chars to ~13.41 chars
#      (Overview columns E & F, and column C on each per-locale sheet).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- 1. Status text: "Ready for handoff" -> "In Translation" ---
$wsOverview.Range("E2:F4").Value = "In Translation"
$wsZhCn.Range("C2:C4").Value     = "In Translation"
$wsDeDe.Range("C2:C4").Value     = "In Translation"

# --- 2. Narrow the status columns ---
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth     = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth     = 12.5
